$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (refresh time 02:05 -> 02:35)
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 02:35"

# Row 4
$ws.Cells.Item(4, 2).Value = 1430348
$ws.Cells.Item(4, 3).Value = 21712
$ws.Cells.Item(4, 4).Value = 310259
$ws.Cells.Item(4, 5).Value = 1034892
$ws.Cells.Item(4, 6).Value = 16349
$ws.Cells.Item(4, 7).Value = 1772
$ws.Cells.Item(4, 8).Value = 85197

# Row 37
$ws.Cells.Item(37, 1).Value = "Japon"
$ws.Cells.Item(37, 2).Value = 16049
$ws.Cells.Item(37, 3).Value = 81
$ws.Cells.Item(37, 4).Value = 8920
$ws.Cells.Item(37, 5).Value = 6451
$ws.Cells.Item(37, 6).Value = 243
$ws.Cells.Item(37, 7).Value = 21
$ws.Cells.Item(37, 8).Value = 678

# Row 38
$ws.Cells.Item(38, 1).Value = "Rumania"
$ws.Cells.Item(38, 2).Value = 16002
$ws.Cells.Item(38, 3).Value = 224
$ws.Cells.Item(38, 4).Value = 7961
$ws.Cells.Item(38, 5).Value = 7005
$ws.Cells.Item(38, 6).Value = 228
$ws.Cells.Item(38, 7).Value = 34
$ws.Cells.Item(38, 8).Value = 1036

# Row 39
$ws.Cells.Item(39, 1).Value = "Austria"
$ws.Cells.Item(39, 2).Value = 15997
$ws.Cells.Item(39, 3).Value = 36
$ws.Cells.Item(39, 4).Value = 14304
$ws.Cells.Item(39, 5).Value = 1069
$ws.Cells.Item(39, 6).Value = 55
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = 624

# Row 100
$ws.Cells.Item(100, 2).Value = 1032
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 759
$ws.Cells.Item(100, 5).Value = 228
$ws.Cells.Item(100, 6).Value = 5
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 45

# Row 101
$ws.Cells.Item(101, 1).Value = "Gabon"
$ws.Cells.Item(101, 2).Value = 1004
$ws.Cells.Item(101, 3).Value = 141
$ws.Cells.Item(101, 4).Value = 152
$ws.Cells.Item(101, 5).Value = 843
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 9

# Row 102
$ws.Cells.Item(102, 1).Value = "Maldivas"
$ws.Cells.Item(102, 2).Value = 955
$ws.Cells.Item(102, 3).Value = 51
$ws.Cells.Item(102, 4).Value = 40
$ws.Cells.Item(102, 5).Value = 911
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(102, 8).Value = 4

# Row 103
$ws.Cells.Item(103, 1).Value = "Letonia"
$ws.Cells.Item(103, 2).Value = 951
$ws.Cells.Item(103, 3).Value = 1
$ws.Cells.Item(103, 4).Value = 627
$ws.Cells.Item(103, 5).Value = 305
$ws.Cells.Item(103, 6).Value = 2
$ws.Cells.Item(103, 7).Value = 1
$ws.Cells.Item(103, 8).Value = 19

# Row 104
$ws.Cells.Item(104, 1).Value = "Republica de Chipre"
$ws.Cells.Item(104, 2).Value = 905
$ws.Cells.Item(104, 3).Value = 2
$ws.Cells.Item(104, 4).Value = 449
$ws.Cells.Item(104, 5).Value = 439
$ws.Cells.Item(104, 6).Value = 10
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 17

# Row 105
$ws.Cells.Item(105, 1).Value = "Sri Lanka"
$ws.Cells.Item(105, 2).Value = 893
$ws.Cells.Item(105, 3).Value = 4
$ws.Cells.Item(105, 4).Value = 382
$ws.Cells.Item(105, 5).Value = 502
$ws.Cells.Item(105, 6).Value = 1
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 9

# Row 106
$ws.Cells.Item(106, 1).Value = "Albania"
$ws.Cells.Item(106, 2).Value = 880
$ws.Cells.Item(106, 3).Value = 4
$ws.Cells.Item(106, 4).Value = 688
$ws.Cells.Item(106, 5).Value = 161
$ws.Cells.Item(106, 6).Value = 1
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 31

# Row 107
$ws.Cells.Item(107, 1).Value = "Libano"
$ws.Cells.Item(107, 2).Value = 878
$ws.Cells.Item(107, 3).Value = 8
$ws.Cells.Item(107, 4).Value = 236
$ws.Cells.Item(107, 5).Value = 616
$ws.Cells.Item(107, 6).Value = 4
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 26

# Row 149
$ws.Cells.Item(149, 1).Value = "Martinica"
$ws.Cells.Item(149, 2).Value = 189
$ws.Cells.Item(149, 3).Value = 2
$ws.Cells.Item(149, 4).Value = 91
$ws.Cells.Item(149, 5).Value = 84
$ws.Cells.Item(149, 6).Value = 2
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 14

# Row 150
$ws.Cells.Item(150, 1).Value = "Suazilandia"
$ws.Cells.Item(150, 2).Value = 187
$ws.Cells.Item(150, 3).Value = 3
$ws.Cells.Item(150, 4).Value = 48
$ws.Cells.Item(150, 5).Value = 137
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 2

# Row 164
$ws.Cells.Item(164, 2).Value = 101
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 91
$ws.Cells.Item(164, 5).Value = 7
$ws.Cells.Item(164, 6).Value = 4
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 3

# Row 192
$ws.Cells.Item(192, 1).Value = "Santa Lucia"
$ws.Cells.Item(192, 2).Value = 18
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 18
$ws.Cells.Item(192, 5).Value = 0
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

# Row 193
$ws.Cells.Item(193, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(193, 2).Value = 18
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 18
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

# Row 194
$ws.Cells.Item(194, 1).Value = "Belice"
$ws.Cells.Item(194, 2).Value = 18
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 16
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 2
